# Scheduled market-data refresh: updates computed price/profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)
# on the affected rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4820.4
$ws.Range("I12").Value = 6666.8335
$ws.Range("J12").Value = 2050.75
$ws.Range("K12").Value = 6666.8335
$ws.Range("L12").Value = 2050.75
$ws.Range("M12").Value = -6496.8335
$ws.Range("N12").Value = -2390.75

$ws.Range("H15").Value = 835.6279
$ws.Range("I15").Value = 835.6279
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2506.8837
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -2337.8837

$ws.Range("H33").Value = 144.23077
$ws.Range("I33").Value = 112.5
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 112.5
$ws.Range("L33").Value = 250
$ws.Range("M33").Value = 116.5
$ws.Range("N33").Value = -708

$ws.Range("H40").Value = 3753.5715
$ws.Range("I40").Value = 1943.75
$ws.Range("J40").Value = 6166.6665
$ws.Range("K40").Value = 1943.75
$ws.Range("L40").Value = 6166.6665
$ws.Range("M40").Value = -1768.75
$ws.Range("N40").Value = -6516.6665

$ws.Range("H51").Value = 8337030
$ws.Range("I51").Value = 3851.7144
$ws.Range("J51").Value = 27781112
$ws.Range("K51").Value = 3851.7144
$ws.Range("L51").Value = 27781112
$ws.Range("M51").Value = -3367.7144
$ws.Range("N51").Value = -27782080

$ws.Range("H74").Value = 5106.4546
$ws.Range("I74").Value = 4726
$ws.Range("J74").Value = 5249.125
$ws.Range("K74").Value = 4726
$ws.Range("L74").Value = 5249.125
$ws.Range("M74").Value = -3790
$ws.Range("N74").Value = -7121.125

$ws.Range("H76").Value = 6149.2
$ws.Range("I76").Value = 6186.5
$ws.Range("J76").Value = 6000
$ws.Range("K76").Value = 6186.5
$ws.Range("L76").Value = 6000
$ws.Range("M76").Value = -5871.5
$ws.Range("N76").Value = -6630

$ws.Range("H77").Value = 5106.4546
$ws.Range("I77").Value = 4726
$ws.Range("J77").Value = 5249.125
$ws.Range("K77").Value = 23630
$ws.Range("L77").Value = 26245.625
$ws.Range("M77").Value = -18950
$ws.Range("N77").Value = -35605.625

$ws.Range("H79").Value = 6149.2
$ws.Range("I79").Value = 6186.5
$ws.Range("J79").Value = 6000
$ws.Range("K79").Value = 6186.5
$ws.Range("L79").Value = 6000
$ws.Range("M79").Value = -5094.5
$ws.Range("N79").Value = -8184

$ws.Range("H107").Value = 188653.75
$ws.Range("I107").Value = 1538.3334
$ws.Range("J107").Value = 750000
$ws.Range("K107").Value = 1538.3334
$ws.Range("L107").Value = 750000
$ws.Range("M107").Value = 381.6666
$ws.Range("N107").Value = -753840

$ws.Range("H135").Value = 78201.38
$ws.Range("I135").Value = 1192.5714
$ws.Range("J135").Value = 168045
$ws.Range("K135").Value = 10733.1426
$ws.Range("L135").Value = 1512405
$ws.Range("M135").Value = -8198.142600000001
$ws.Range("N135").Value = -1517475

$ws.Range("H138").Value = 2507.8096
$ws.Range("I138").Value = 3295.0833
$ws.Range("J138").Value = 2192.9
$ws.Range("K138").Value = 9885.249899999999
$ws.Range("L138").Value = 6578.700000000001
$ws.Range("M138").Value = -4745.249899999999
$ws.Range("N138").Value = -16858.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6129.6206
$ws.Range("I45").Value = 8802.933999999999
$ws.Range("J45").Value = 3265.3572
$ws.Range("K45").Value = 8802.933999999999
$ws.Range("L45").Value = 3265.3572
$ws.Range("M45").Value = -8425.933999999999
$ws.Range("N45").Value = -4019.3572

$ws.Range("H61").Value = 2898.577
$ws.Range("I61").Value = 2528.2354
$ws.Range("J61").Value = 3598.111
$ws.Range("K61").Value = 2528.2354
$ws.Range("L61").Value = 3598.111
$ws.Range("M61").Value = -2316.2354
$ws.Range("N61").Value = -4022.111

$ws.Range("H74").Value = 2053.25
$ws.Range("I74").Value = 1830.4634
$ws.Range("J74").Value = 3358.1428
$ws.Range("K74").Value = 1830.4634
$ws.Range("L74").Value = 3358.1428
$ws.Range("M74").Value = -956.4634000000001
$ws.Range("N74").Value = -5106.1428

$ws.Range("H77").Value = 2053.25
$ws.Range("I77").Value = 1830.4634
$ws.Range("J77").Value = 3358.1428
$ws.Range("K77").Value = 9152.317000000001
$ws.Range("L77").Value = 16790.714
$ws.Range("M77").Value = -4784.317000000001
$ws.Range("N77").Value = -25526.714

$ws.Range("H122").Value = 3064
$ws.Range("I122").Value = 3082.875
$ws.Range("J122").Value = 2913
$ws.Range("K122").Value = 9248.625
$ws.Range("L122").Value = 8739
$ws.Range("M122").Value = -6798.625
$ws.Range("N122").Value = -13639

$ws.Range("H132").Value = 3183.08
$ws.Range("I132").Value = 1254.8334
$ws.Range("J132").Value = 8141.4287
$ws.Range("K132").Value = 3764.5002
$ws.Range("L132").Value = 24424.2861
$ws.Range("M132").Value = -1234.5002
$ws.Range("N132").Value = -29484.2861

$ws.Range("H136").Value = 2898.577
$ws.Range("I136").Value = 2528.2354
$ws.Range("J136").Value = 3598.111
$ws.Range("K136").Value = 7584.706200000001
$ws.Range("L136").Value = 10794.333
$ws.Range("M136").Value = -5034.706200000001
$ws.Range("N136").Value = -15894.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16436
$ws.Range("I20").Value = 14275.444
$ws.Range("J20").Value = 21297.25
$ws.Range("K20").Value = 14275.444
$ws.Range("L20").Value = 21297.25
$ws.Range("M20").Value = -14028.444
$ws.Range("N20").Value = -21791.25

$ws.Range("H105").Value = 1395.174
$ws.Range("I105").Value = 1218.5333
$ws.Range("J105").Value = 1726.375
$ws.Range("K105").Value = 1218.5333
$ws.Range("L105").Value = 1726.375
$ws.Range("M105").Value = 528.4666999999999
$ws.Range("N105").Value = -5220.375

$ws.Range("H134").Value = 1105.7222
$ws.Range("I134").Value = 955.93335
$ws.Range("J134").Value = 1854.6666
$ws.Range("K134").Value = 2867.80005
$ws.Range("L134").Value = 5563.9998
$ws.Range("M134").Value = -332.8000499999998
$ws.Range("N134").Value = -10633.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1117.6
$ws.Range("I6").Value = 899.5
$ws.Range("J6").Value = 1990
$ws.Range("K6").Value = 899.5
$ws.Range("L6").Value = 1990
$ws.Range("M6").Value = -786.5
$ws.Range("N6").Value = -2216

$ws.Range("H16").Value = 1969.1875
$ws.Range("I16").Value = 1900.5
$ws.Range("J16").Value = 2450
$ws.Range("K16").Value = 1900.5
$ws.Range("L16").Value = 2450
$ws.Range("M16").Value = -1613.5
$ws.Range("N16").Value = -3024

$ws.Range("H62").Value = 17032
$ws.Range("I62").Value = 18660.25
$ws.Range("J62").Value = 4006
$ws.Range("K62").Value = 18660.25
$ws.Range("L62").Value = 4006
$ws.Range("M62").Value = -18036.25
$ws.Range("N62").Value = -5254

$ws.Range("H65").Value = 17032
$ws.Range("I65").Value = 18660.25
$ws.Range("J65").Value = 4006
$ws.Range("K65").Value = 93301.25
$ws.Range("L65").Value = 20030
$ws.Range("M65").Value = -90181.25
$ws.Range("N65").Value = -26270

$ws.Range("H94").Value = 1112.1428
$ws.Range("I94").Value = 207.75
$ws.Range("J94").Value = 1473.9
$ws.Range("K94").Value = 207.75
$ws.Range("L94").Value = 1473.9
$ws.Range("M94").Value = 243.25
$ws.Range("N94").Value = -2375.9

$ws.Range("H105").Value = 2698
$ws.Range("I105").Value = 2726.5715
$ws.Range("J105").Value = 2498
$ws.Range("K105").Value = 2726.5715
$ws.Range("L105").Value = 2498
$ws.Range("M105").Value = -979.5715
$ws.Range("N105").Value = -5992

$ws.Range("H107").Value = 1717.75
$ws.Range("I107").Value = 1426.5264
$ws.Range("J107").Value = 2332.5557
$ws.Range("K107").Value = 1426.5264
$ws.Range("L107").Value = 2332.5557
$ws.Range("M107").Value = 493.4736
$ws.Range("N107").Value = -6172.5557

$ws.Range("H113").Value = 1969.1875
$ws.Range("I113").Value = 1900.5
$ws.Range("J113").Value = 2450
$ws.Range("K113").Value = 1900.5
$ws.Range("L113").Value = 2450
$ws.Range("M113").Value = 269.5
$ws.Range("N113").Value = -6790

$ws.Range("H132").Value = 2695.9048
$ws.Range("I132").Value = 2673.4062
$ws.Range("J132").Value = 2767.9
$ws.Range("K132").Value = 8020.2186
$ws.Range("L132").Value = 8303.700000000001
$ws.Range("M132").Value = -5490.2186
$ws.Range("N132").Value = -13363.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 963.7826
$ws.Range("I107").Value = 693.4
$ws.Range("J107").Value = 1038.8889
$ws.Range("K107").Value = 2080.2
$ws.Range("L107").Value = 3116.6667
$ws.Range("M107").Value = -160.1999999999998
$ws.Range("N107").Value = -6956.6667

$ws.Range("H113").Value = 1703.3846
$ws.Range("I113").Value = 1150
$ws.Range("J113").Value = 1949.3334
$ws.Range("K113").Value = 3450
$ws.Range("L113").Value = 5848.0002
$ws.Range("M113").Value = -1280
$ws.Range("N113").Value = -10188.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3334408.2
$ws.Range("I3").Value = 3333483
$ws.Range("J3").Value = 3335333.2
$ws.Range("K3").Value = 3333483
$ws.Range("L3").Value = 3335333.2
$ws.Range("M3").Value = -3333367
$ws.Range("N3").Value = -3335565.2

$ws.Range("H80").Value = 3521.625
$ws.Range("I80").Value = 2833.7778
$ws.Range("J80").Value = 4406
$ws.Range("K80").Value = 2833.7778
$ws.Range("L80").Value = 4406
$ws.Range("M80").Value = -1835.7778
$ws.Range("N80").Value = -6402

$ws.Range("H83").Value = 3521.625
$ws.Range("I83").Value = 2833.7778
$ws.Range("J83").Value = 4406
$ws.Range("K83").Value = 14168.889
$ws.Range("L83").Value = 22030
$ws.Range("M83").Value = -9176.888999999999
$ws.Range("N83").Value = -32014

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 37505
$ws.Range("I3").Value = 25004
$ws.Range("J3").Value = 50006
$ws.Range("K3").Value = 25004
$ws.Range("L3").Value = 50006
$ws.Range("M3").Value = -24892
$ws.Range("N3").Value = -50230

$ws.Range("H7").Value = 6509.6343
$ws.Range("I7").Value = 8274.875
$ws.Range("J7").Value = 5379.88
$ws.Range("K7").Value = 8274.875
$ws.Range("L7").Value = 5379.88
$ws.Range("M7").Value = -8162.875
$ws.Range("N7").Value = -5603.88

$ws.Range("H15").Value = 37505
$ws.Range("I15").Value = 25004
$ws.Range("J15").Value = 50006
$ws.Range("K15").Value = 25004
$ws.Range("L15").Value = 50006
$ws.Range("M15").Value = -24834
$ws.Range("N15").Value = -50346

$ws.Range("H16").Value = 772.8889
$ws.Range("I16").Value = 774.8570999999999
$ws.Range("J16").Value = 766
$ws.Range("K16").Value = 774.8570999999999
$ws.Range("L16").Value = 766
$ws.Range("M16").Value = -604.8570999999999
$ws.Range("N16").Value = -1106

$ws.Range("H23").Value = 11000
$ws.Range("I23").Value = 11000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 11000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -10770

$ws.Range("H55").Value = 247.35294
$ws.Range("I55").Value = 254.15384
$ws.Range("J55").Value = 225.25
$ws.Range("K55").Value = 254.15384
$ws.Range("L55").Value = 225.25
$ws.Range("M55").Value = -81.15384
$ws.Range("N55").Value = -571.25

$ws.Range("H74").Value = 20108.5
$ws.Range("I74").Value = 15000
$ws.Range("J74").Value = 25217
$ws.Range("K74").Value = 15000
$ws.Range("L74").Value = 25217
$ws.Range("M74").Value = -14002
$ws.Range("N74").Value = -27213

$ws.Range("H77").Value = 20108.5
$ws.Range("I77").Value = 15000
$ws.Range("J77").Value = 25217
$ws.Range("K77").Value = 45000
$ws.Range("L77").Value = 75651
$ws.Range("M77").Value = -40008
$ws.Range("N77").Value = -85635

$ws.Range("H126").Value = 6509.6343
$ws.Range("I126").Value = 8274.875
$ws.Range("J126").Value = 5379.88
$ws.Range("K126").Value = 24824.625
$ws.Range("L126").Value = 16139.64
$ws.Range("M126").Value = -22354.625
$ws.Range("N126").Value = -21079.64

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 7266298.5
$ws.Range("I3").Value = 12614797
$ws.Range("J3").Value = 134966.67
$ws.Range("K3").Value = 12614797
$ws.Range("L3").Value = 134966.67
$ws.Range("M3").Value = -12614683
$ws.Range("N3").Value = -135194.67

$ws.Range("H132").Value = 4461.8857
$ws.Range("I132").Value = 4550.3335
$ws.Range("J132").Value = 3002.5
$ws.Range("K132").Value = 13651.0005
$ws.Range("L132").Value = 9007.5
$ws.Range("M132").Value = -11121.0005
$ws.Range("N132").Value = -14067.5
